$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.645.87"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "1.874.01"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -0.87%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.48"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5078"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3921"
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08381"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.29"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.107"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.195"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "1.877.01"
$ws.Range("E13").Value = "  +3.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.36"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.264"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.12"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001100"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06714"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.64"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.932"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "28.670.35"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.09"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.191"
$ws.Range("E25").Value = "  -3.78%  "
$ws.Range("D26").Value = "2.084.85"
$ws.Range("E26").Value = "  +3.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.33"
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.69"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.423"
$ws.Range("E29").Value = "  +3.04%  "
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1037"
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.046"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.780"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.623"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02455"
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.017"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2163"
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.041"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.190"
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.242"
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6383"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5992"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.05"
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.674"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.005"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.220"
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.196"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.27"
$ws.Range("E51").Value = "  +1.38%  "
